$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 3): "kimppakohde" (collective/partnership site) test fixture ---
$ws.Range("A3").Value = "0000000-9"
$ws.Range("E3").Value = "01-0000123-01"
$ws.Range("F3").Value = "134567890B"
$ws.Range("G3").Value = "KIRKKOÄYRÄÄNTIE 1D"
$ws.Range("H3").Value = "16200 ARTJÄRVI"
$ws.Range("K3").Value = "JOHN LINDROTH"
$ws.Range("L3").Value = "JOHN LINDROTH"
$ws.Range("M3").Value = "KIRKKOÄYRÄÄNTIE 1D"
$ws.Range("N3").Value = "16200 ARTJÄRVI"
$ws.Range("O3").Value = "FI"

# Date-like text values must be forced to text so they are stored as literal
# strings (matching the source row) instead of being parsed into date serials.
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "1.1.2023"
$ws.Range("Q3").NumberFormat = "General"

$ws.Range("R3").NumberFormat = "@"
$ws.Range("R3").Value = "31.12.2023"
$ws.Range("R3").NumberFormat = "General"

$ws.Range("S3").Value = "Energia"
$ws.Range("T3").Value = "660 L ENERGIAJÄTEASTIA TYHJENNYS"
$ws.Range("U3").Value = 66
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = "0,66"
$ws.Range("X3").Value = 8
$ws.Range("Y3").Value = 1
$ws.Range("AA3").Value = 1
$ws.Range("AB3").Value = 53

$ws.Range("AC3").NumberFormat = "@"
$ws.Range("AC3").Value = "1.1.1900"
$ws.Range("AC3").NumberFormat = "General"

$ws.Range("AD3").NumberFormat = "@"
$ws.Range("AD3").Value = "31.12.2099"
$ws.Range("AD3").NumberFormat = "General"

$ws.Range("AK3").Value = "01-0000999-99"
$ws.Range("AL3").Value = "KIRKKOÄYRÄÄN ENERGIAKIMPPA"
$ws.Range("AM3").Value = "01-1110999-01"
$ws.Range("AN3").Value = "KIRKKOÄYRÄÄNTIE 11"
$ws.Range("AO3").Value = "16200 ARTJÄRVI"
$ws.Range("AP3").Value = 560

# --- Column width adjustments (new/wider columns to fit the new data) ---
$ws.Columns.Item(7).ColumnWidth = 25
$ws.Columns.Item(20).ColumnWidth = 44.83
$ws.Columns.Item(14).ColumnWidth = 22.25
$ws.Columns.Item(37).ColumnWidth = 16.55
$ws.Columns.Item(38).ColumnWidth = 33.11
$ws.Columns.Item(39).ColumnWidth = 16.55
$ws.Columns.Item(40).ColumnWidth = 21.84
$ws.Columns.Item(41).ColumnWidth = 22.39

# --- Selection moves to the next empty row, as it would after data entry ---
$ws.Range("A4").Select()
